$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''34.733.81'
$ws.Range('E2').Value = '  +0.84%  '

$ws.Range('D3').Value = '''1.811.79'
$ws.Range('E3').Value = '  +0.39%  '

$ws.Range('D4').Value = '''0.997'
$ws.Range('E4').Value = '  -0.72%  '

$ws.Range('D5').Value = '''224.99'
$ws.Range('E5').Value = '  -1.37%  '

$ws.Range('D6').Value = '''0.607'
$ws.Range('E6').Value = '  +0.94%  '

$ws.Range('D7').Value = '''0.998'
$ws.Range('E7').Value = '  -0.66%  '

$ws.Range('D8').Value = '''39.41'
$ws.Range('E8').Value = '  +8.73%  '

$ws.Range('D9').Value = '''0.291'
$ws.Range('E9').Value = '  -2.51%  '

$ws.Range('D10').Value = '''0.0672'
$ws.Range('E10').Value = '  -3.66%  '

$ws.Range('D11').Value = '''0.100'
$ws.Range('E11').Value = '  +3.59%  '

$ws.Range('D12').Value = '''2.062.76'
$ws.Range('E12').Value = '  -0.25%  '

$ws.Range('D13').Value = '''1.800.23'
$ws.Range('E13').Value = '  -0.68%  '

$ws.Range('D14').Value = '''10.99'
$ws.Range('E14').Value = '  -2.93%  '

$ws.Range('D15').Value = '''0.635'
$ws.Range('E15').Value = '  -2.05%  '

$ws.Range('D16').Value = '''34.658.98'
$ws.Range('E16').Value = '  +0.58%  '

$ws.Range('D17').Value = '''4.40'
$ws.Range('E17').Value = '  -1.48%  '

$ws.Range('D18').Value = '''68.31'
$ws.Range('E18').Value = '  -2.09%  '

$ws.Range('D19').Value = '''241.79'
$ws.Range('E19').Value = '  -1.71%  '

$ws.Range('D20').Value = '''0.0₃0769'
$ws.Range('E20').Value = '  -2.68%  '

$ws.Range('D21').Value = '''11.14'
$ws.Range('E21').Value = '  -3.26%  '

$ws.Range('D22').Value = '''0.999'
$ws.Range('E22').Value = '  -0.51%  '

$ws.Range('D23').Value = '''4.11'
$ws.Range('E23').Value = '  -1.65%  '

$ws.Range('D24').Value = '''2.18'
$ws.Range('E24').Value = '  -1.96%  '

$ws.Range('D25').Value = '''171.08'
$ws.Range('E25').Value = '  -0.07%  '

$ws.Range('D26').Value = '''7.72'
$ws.Range('E26').Value = '  -4.51%  '

$ws.Range('D27').Value = '''17.63'
$ws.Range('E27').Value = '  +0.64%  '

$ws.Range('D28').Value = '''0.122'
$ws.Range('E28').Value = '  +0.35%  '

$ws.Range('D29').Value = '''0.998'
$ws.Range('E29').Value = '  -0.70%  '

$ws.Range('D30').Value = '''1.23'
$ws.Range('E30').Value = '  -1.70%  '

$ws.Range('D31').Value = '''3.76'
$ws.Range('E31').Value = '  -1.82%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '''3.86'
$ws.Range('E32').Value = '  -3.27%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '''0.0514'
$ws.Range('E33').Value = '  -2.15%  '

$ws.Range('D34').Value = '''1.82'
$ws.Range('E34').Value = '  +1.38%  '

$ws.Range('D35').Value = '''0.646'
$ws.Range('E35').Value = '  -2.26%  '

$ws.Range('D36').Value = '''1.315.10'
$ws.Range('E36').Value = '  -5.11%  '

$ws.Range('D37').Value = '''1.06'
$ws.Range('E37').Value = '  +0.09%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '''2.35'
$ws.Range('E38').Value = '  -1.42%  '

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.0188'
$ws.Range('E39').Value = '  -0.01%  '

$ws.Range('E40').Value = '  +3.44%  '

$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '''82.64'
$ws.Range('E41').Value = '  +0.50%  '

$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').Value = '''2.44'
$ws.Range('E42').Value = '  +0.58%  '

$ws.Range('D43').Value = '''2.81'
$ws.Range('E43').Value = '  -1.33%  '

$ws.Range('D44').Value = '''0.950'
$ws.Range('E44').Value = '  -0.41%  '

$ws.Range('D45').Value = '''14.20'
$ws.Range('E45').Value = '  +5.21%  '

$ws.Range('E46').Value = '  +3.62%  '

$ws.Range('D47').Value = '''1.966.55'
$ws.Range('E47').Value = '  -0.12%  '

$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').Value = '''5.72'
$ws.Range('E48').Value = '  -4.38%  '

$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').Value = '''0.999'
$ws.Range('E49').Value = '  -0.56%  '

$ws.Range('D50').Value = '''102.08'
$ws.Range('E50').Value = '  -1.01%  '

$ws.Range('D51').Value = '''0.0613'
$ws.Range('E51').Value = '  -0.21%  '
